$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '30.602.46'
Set-TextValue 'E2' '  +0.37%  '
Set-TextValue 'D3' '1.881.63'
Set-TextValue 'E3' '  +0.10%  '
Set-TextValue 'E4' '  +0.07%  '
Set-TextValue 'D5' '249.71'
Set-TextValue 'E5' '  +1.29%  '
Set-TextValue 'D6' '1.000'
Set-TextValue 'E6' '  +0.07%  '
Set-TextValue 'E7' '  -0.17%  '
Set-TextValue 'D8' '0.2932'
Set-TextValue 'E8' '  +0.99%  '
Set-TextValue 'D9' '0.06525'
Set-TextValue 'D10' '21.90'
Set-TextValue 'E10' '  -0.22%  '
Set-TextValue 'D11' '0.07747'
Set-TextValue 'E11' '  +0.02%  '
Set-TextValue 'D12' '96.95'
Set-TextValue 'E12' '  -0.09%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.882.21'
Set-TextValue 'E13' '  +0.23%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D14' '0.7374'
Set-TextValue 'E14' '  -0.83%  '
Set-TextValue 'E15' '  +2.38%  '
Set-TextValue 'E16' '  +0.10%  '
Set-TextValue 'D17' '30.727.14'
Set-TextValue 'E17' '  +0.83%  '
Set-TextValue 'E18' '  -3.44%  '
Set-TextValue 'D19' '0.000007532'
Set-TextValue 'E19' '  -0.65%  '
Set-TextValue 'D20' '1.000'
Set-TextValue 'E20' '  +0.01%  '
Set-TextValue 'D21' '2.129.37'
Set-TextValue 'E21' '  +0.30%  '
Set-TextValue 'D22' '5.313'
Set-TextValue 'D23' '1.001'
Set-TextValue 'E23' '  +0.09%  '
Set-TextValue 'D24' '6.231'
Set-TextValue 'E24' '  +0.93%  '
Set-TextValue 'D25' '9.203'
Set-TextValue 'E25' '  -0.92%  '
Set-TextValue 'D26' '163.91'
Set-TextValue 'E26' '  -0.35%  '
Set-TextValue 'E27' '  -0.31%  '
Set-TextValue 'D28' '1.909'
Set-TextValue 'E28' '  -2.64%  '
Set-TextValue 'D29' '1.342'
Set-TextValue 'E29' '  -2.14%  '
Set-TextValue 'D30' '0.09697'
Set-TextValue 'E30' '  -2.95%  '
Set-TextValue 'D31' '1.510'
Set-TextValue 'E31' '  -0.26%  '
Set-TextValue 'D32' '4.286'
Set-TextValue 'E32' '  -1.10%  '
Set-TextValue 'D33' '4.141'
Set-TextValue 'E33' '  +1.97%  '
Set-TextValue 'D34' '0.04860'
Set-TextValue 'E34' '  +1.58%  '
Set-TextValue 'D35' '1.127'
Set-TextValue 'E35' '  +0.07%  '
Set-TextValue 'D36' '0.6987'
Set-TextValue 'E36' '  +0.00%  '
Set-TextValue 'E37' '  +0.19%  '
Set-TextValue 'D38' '0.01904'
Set-TextValue 'E38' '  +1.92%  '
Set-TextValue 'D39' '2.780'
Set-TextValue 'E39' '  +1.84%  '
Set-TextValue 'D40' '6.323'
Set-TextValue 'E40' '  -0.52%  '
Set-TextValue 'D41' '74.94'
Set-TextValue 'E41' '  +6.93%  '
Set-TextValue 'D42' '2.012'
Set-TextValue 'E42' '  +3.85%  '
Set-TextValue 'D43' '0.4239'
Set-TextValue 'E43' '  +1.45%  '
Set-TextValue 'D44' '0.8417'
Set-TextValue 'E44' '  +0.42%  '
Set-TextValue 'E45' '  +0.05%  '
Set-TextValue 'D46' '102.37'
Set-TextValue 'E46' '  -0.49%  '
Set-TextValue 'D47' '9.400'
Set-TextValue 'E47' '  +0.68%  '
Set-TextValue 'D48' '7.053'
Set-TextValue 'E48' '  -0.63%  '
Set-TextValue 'D49' '35.63'
Set-TextValue 'E49' '  +0.72%  '
Set-TextValue 'D50' '913.96'
Set-TextValue 'E50' '  -0.88%  '
Set-TextValue 'E51' '  +2.26%  '
